# Adds documentation of the notify() Artisan Command to the Commands sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Insert a new row above row 85 (shifting existing row 85+ down by one row)
$ws.Rows.Item(85).Insert()

$ws.Cells.Item(85, 2).Value = "notify(<title>,[<msg>])"
$ws.Cells.Item(85, 3).Value = "sends notification with title <title> and optional message <msg>"

$ws.Range("B85:C85").Select()
